$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 161

for ($r = 1; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)

    $aVal = $aCell.Value2
    $bVal = $bCell.Value2

    $aCell.Value = $bVal
    $bCell.Value = $aVal

    if ($r -gt 1) {
        $eCell = $ws.Cells.Item($r, 5)
        $eVal = $eCell.Text
        if ($eVal -ne $null -and $eVal -ne "") {
            $eCell.Value = "$eVal m²"
        }
    }
}
